$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MySingleFunction")
Write-Host ($ws.Columns.Item(13) | Get-Member | Out-String) 
